$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Priority column (A) updates
$ws.Range("A20").Value = 9
$ws.Range("A21").Value = 9
$ws.Range("A57").Value = 1
$ws.Range("A59").Value = 8

# comments column (E) text updates
$ws.Range("E47").Value = "out of date. also, when is it null vs when is it zero? see uid CYP19A1 vs uid CYP19A1_NC in pop_stats"
$ws.Range("E57").Value = "should be simple. unclear how CRAVAT decides when to make a mupit link.  "

# Update selected/active cell to match new view state
$ws.Range("C49").Select()
